$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3097678110427751
$ws.Range("C2").Value = 0.05372333215592562
$ws.Range("E2").Value = 0.41800318993063
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.002411376888213563
$ws.Range("I2").Value = 0.4307230857009898
$ws.Range("K2").Value = 0.3267177242355785
$ws.Range("N2").Value = 1.179559757958838
$ws.Range("O2").Value = 1.91791305151223
$ws.Range("B3").Value = 0.2727858857300021
$ws.Range("C3").Value = 0.04746014237262841
$ws.Range("E3").Value = 0.3647724438965838
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.002413605453381719
$ws.Range("I3").Value = 0.4357572175710089
$ws.Range("K3").Value = 0.2856613219461792
$ws.Range("N3").Value = 1.192673224353818
$ws.Range("O3").Value = 1.931545568098173
$ws.Range("B4").Value = 0.2500647799018054
$ws.Range("C4").Value = 0.04359434711180654
$ws.Range("E4").Value = 0.3321659454530561
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.00241504506506099
$ws.Range("I4").Value = 0.4391146096214591
$ws.Range("K4").Value = 0.2604078710094768
$ws.Range("N4").Value = 1.201175816187089
$ws.Range("O4").Value = 1.94121724142903
$ws.Range("B5").Value = 0.2408026504438681
$ws.Range("C5").Value = 0.04201394151434101
$ws.Range("E5").Value = 0.3188965596776683
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002415649694157851
$ws.Range("I5").Value = 0.4405496746317539
$ws.Range("K5").Value = 0.2501060494026603
$ws.Range("N5").Value = 1.204754071335945
$ws.Range("O5").Value = 1.945485154583992
$ws.Range("B6").Value = 0.2392645069968466
$ws.Range("C6").Value = 0.04175121089608069
$ws.Range("E6").Value = 0.316694235281048
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002415751179725993
$ws.Range("I6").Value = 0.4407920040331206
$ws.Range("K6").Value = 0.2483947969808327
$ws.Range("N6").Value = 1.205355085947922
$ws.Range("O6").Value = 1.946213548820523
$ws.Range("B7").Value = 0.249939879529876
$ws.Range("C7").Value = 0.04357305369393316
$ws.Range("E7").Value = 0.3319869187683651
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.002415053146602443
$ws.Range("I7").Value = 0.4391336926551901
$ws.Range("K7").Value = 0.260268980389526
$ws.Range("N7").Value = 1.201223614711353
$ws.Range("O7").Value = 1.941273478274283
$ws.Range("B8").Value = 0.2970196122463733
$ws.Range("C8").Value = 0.05156797474487007
$ws.Range("E8").Value = 0.3996322856658026
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.002412130540255862
$ws.Range("I8").Value = 0.4324035049595629
$ws.Range("K8").Value = 0.3125709678038504
$ws.Range("N8").Value = 1.183987706509555
$ws.Range("O8").Value = 1.922343114778386
$ws.Range("B9").Value = 0.389216527157231
$ws.Range("C9").Value = 0.06708627123121857
$ws.Range("E9").Value = 0.5329696169925171
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.002406962212590887
$ws.Range("I9").Value = 0.4213236290494926
$ws.Range("K9").Value = 0.4147684117178869
$ws.Range("N9").Value = 1.153765037388027
$ws.Range("O9").Value = 1.895573103188553
$ws.Range("B10").Value = 0.4568634257641975
$ws.Range("C10").Value = 0.07839175826399014
$ws.Range("E10").Value = 0.6314571581954311
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002403504611052366
$ws.Range("I10").Value = 0.414480119320455
$ws.Range("K10").Value = 0.4896202519636574
$ws.Range("N10").Value = 1.133739751667484
$ws.Range("O10").Value = 1.882256113145161
$ws.Range("B11").Value = 0.4876158002136037
$ws.Range("C11").Value = 0.08351452198242271
$ws.Range("E11").Value = 0.6763992865148651
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.00240200462779828
$ws.Range("I11").Value = 0.4116496728655221
$ws.Range("K11").Value = 0.5236203150533925
$ws.Range("N11").Value = 1.125102722915369
$ws.Range("O11").Value = 1.877585682767602
$ws.Range("B12").Value = 0.4992576198761753
$ws.Range("C12").Value = 0.0854514892624394
$ws.Range("E12").Value = 0.6934395173588825
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002401447048140348
$ws.Range("I12").Value = 0.4106186151233118
$ws.Range("K12").Value = 0.5364877133582695
$ws.Range("N12").Value = 1.121900072754674
$ws.Range("O12").Value = 1.87601733701581
$ws.Range("B13").Value = 0.4967505062640782
$ws.Range("C13").Value = 0.08503445842788437
$ws.Range("E13").Value = 0.6897686124494697
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002401566669665315
$ws.Range("I13").Value = 0.4108388566590548
$ws.Range("K13").Value = 0.5337168377752732
$ws.Range("N13").Value = 1.122586795394188
$ws.Range("O13").Value = 1.876346190902211
$ws.Range("B14").Value = 0.4885736513741676
$ws.Range("C14").Value = 0.08367393588342509
$ws.Range("E14").Value = 0.6778007543532709
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002401958546701612
$ws.Range("I14").Value = 0.4115640293317924
$ws.Range("K14").Value = 0.5246790802323744
$ws.Range("N14").Value = 1.124837875738098
$ws.Range("O14").Value = 1.877452636682989
$ws.Range("B15").Value = 0.4835646297084679
$ws.Range("C15").Value = 0.08284019690651689
$ws.Range("E15").Value = 0.6704729574810528
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002402199939539621
$ws.Range("I15").Value = 0.412013531646803
$ws.Range("K15").Value = 0.519142178640152
$ws.Range("N15").Value = 1.126225585258918
$ws.Range("O15").Value = 1.878156465862105
$ws.Range("B16").Value = 0.4548532310957398
$ws.Range("C16").Value = 0.07805656680753259
$ws.Range("E16").Value = 0.6285230259818633
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.00240360410073442
$ws.Range("I16").Value = 0.4146707944610455
$ws.Range("K16").Value = 0.4873972159168716
$ws.Range("N16").Value = 1.134313715369657
$ws.Range("O16").Value = 1.882589318892173
$ws.Range("B17").Value = 0.4372341128662356
$ws.Range("C17").Value = 0.07511679184311504
$ws.Range("E17").Value = 0.6028249622794988
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.002404484142631002
$ws.Range("I17").Value = 0.4163734299290276
$ws.Range("K17").Value = 0.4679094848360705
$ws.Range("N17").Value = 1.139396592998729
$ws.Range("O17").Value = 1.885664622197254
$ws.Range("B18").Value = 0.4270981576518125
$ws.Range("C18").Value = 0.07342401434047474
$ws.Range("E18").Value = 0.5880571676569986
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002404997183612171
$ws.Range("I18").Value = 0.4173793427696566
$ws.Range("K18").Value = 0.4566959381816673
$ws.Range("N18").Value = 1.1423646074945
$ws.Range("O18").Value = 1.88756401346879
$ws.Range("B19").Value = 0.4236659850529634
$ws.Range("C19").Value = 0.07285054361749133
$ws.Range("E19").Value = 0.5830592410652997
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.0024051720710453
$ws.Range("I19").Value = 0.4177244928726935
$ws.Range("K19").Value = 0.4528984268178249
$ws.Range("N19").Value = 1.143377163846331
$ws.Range("O19").Value = 1.888229516366039
$ws.Range("B20").Value = 0.4391098986885424
$ws.Range("C20").Value = 0.07542993230534023
$ws.Range("E20").Value = 0.6055592030983092
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.002404389750487702
$ws.Range("I20").Value = 0.4161894274131939
$ws.Range("K20").Value = 0.4699844786741494
$ws.Range("N20").Value = 1.138850908368529
$ws.Range("O20").Value = 1.885323733425565
$ws.Range("B21").Value = 0.4909754889945646
$ws.Range("C21").Value = 0.0840736334732668
$ws.Range("E21").Value = 0.6813154052249786
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.002401843160247965
$ws.Range("I21").Value = 0.4113499210780702
$ws.Range("K21").Value = 0.5273339008257665
$ws.Range("N21").Value = 1.12417483282043
$ws.Range("O21").Value = 1.877122206400372
$ws.Range("B22").Value = 0.524852392949839
$ws.Range("C22").Value = 0.08970581173363712
$ws.Range("E22").Value = 0.7309532169787474
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002400239594024163
$ws.Range("I22").Value = 0.4084247235444529
$ws.Range("K22").Value = 0.5647700516634018
$ws.Range("N22").Value = 1.114979580718654
$ws.Range("O22").Value = 1.872929550886141
$ws.Range("B23").Value = 0.5067736647141317
$ws.Range("C23").Value = 0.0867013701528947
$ws.Range("E23").Value = 0.7044484724931408
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.002401089903044502
$ws.Range("I23").Value = 0.4099641675476775
$ws.Range("K23").Value = 0.544793936299925
$ws.Range("N23").Value = 1.119850970861492
$ws.Range("O23").Value = 1.875060183742221
$ws.Range("B24").Value = 0.4382618760957939
$ws.Range("C24").Value = 0.07528836979298603
$ws.Range("E24").Value = 0.6043230331439418
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002404432403073454
$ws.Range("I24").Value = 0.4162725307133321
$ws.Range("K24").Value = 0.4690464044246028
$ws.Range("N24").Value = 1.139097469721605
$ws.Range("O24").Value = 1.885477440105802
$ws.Range("B25").Value = 0.3642897083277319
$ws.Range("C25").Value = 0.0629050999815064
$ws.Range("E25").Value = 0.4968139809892165
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.002408300492604928
$ws.Range("I25").Value = 0.4240936972118057
$ws.Range("K25").Value = 0.3871614929263387
$ws.Range("N25").Value = 1.161558282986704
$ws.Range("O25").Value = 1.901702702176976
